$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp column (A) and hora_actual column (B) for rows 2-4
$ws.Range("A2:A4").Value = "09/01/2026 14:55:29"
$ws.Range("B2:B4").Value = "14:55"

# Update hora_eta column (C) for rows 2-4 individually
$ws.Range("C2").Value = "15:00"
$ws.Range("C3").Value = "15:10"
$ws.Range("C4").Value = "15:21"
